{"js": "// Find the paragraph that holds the GitHub repository URL and give it the\n// same \"underline + blue font\" treatment that was applied in the commit:\n//   - the paragraph mark itself becomes underlined (selecting the whole\n//     line, pilcrow included, and pressing Ctrl+U)\n//   - the run carrying the URL text becomes underlined AND blue (0070C0)\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst urlNeedle = \"github.com/Gelda3273/Tarea-1\";\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.indexOf(urlNeedle) !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // 1) Underline the paragraph (this also marks the paragraph-mark run\n  //    properties, matching <w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr>).\n  targetParagraph.font.underline = Word.UnderlineType.single;\n  await context.sync();\n\n  // 2) Color just the visible text blue. Using a text search keeps the\n  //    matched range limited to the run itself (no paragraph mark), so the\n  //    paragraph-mark run properties stay underline-only, exactly like the\n  //    diff shows.\n  const searchResults = body.search(\n    \"https://github.com/Gelda3273/Tarea-1/blob/main/Pilas%20est%C3%A1ticas%20C.docx\",\n    { matchCase: false }\n  );\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length > 0) {\n    searchResults.items[0].font.color = \"#0070C0\";\n  } else {\n    // Fallback: color the whole paragraph's range if the exact text wasn't found.\n    targetParagraph.getRange().font.color = \"#0070C0\";\n  }\n  await context.sync();\n}\n", "ps1": "# Find the paragraph that holds the GitHub repository URL and give it the\n# same \"underline + blue font\" treatment that was applied in the commit:\n#   - the paragraph mark itself becomes underlined (as if the whole line,\n#     pilcrow included, had been selected and Ctrl+U pressed)\n#   - the run carrying the URL text becomes underlined AND blue (0070C0)\n\n$d = $word.ActiveDocument\n\n$needle = \"github.com/Gelda3273/Tarea-1\"\n$urlText = \"https://github.com/Gelda3273/Tarea-1/blob/main/Pilas%20est%C3%A1ticas%20C.docx\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$needle*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # 1) Underline the whole paragraph (pilcrow included) -> writes\n    #    <w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr> and underlines\n    #    the run itself too.\n    $target.Range.Font.Underline = [Microsoft.Office.Interop.Word.WdUnderline]::wdUnderlineSingle\n\n    # 2) Color just the URL text blue (0070C0). Scope the Find to a\n    #    duplicate of the paragraph's own range so the paragraph mark is\n    #    left alone (only the run keeps the color).\n    $textRange = $target.Range.Duplicate()\n    $found = $textRange.Find.Execute($urlText)\n    if ($found) {\n        $textRange.Font.Color = 0xC07000\n    } else {\n        $target.Range.Font.Color = 0xC07000\n    }\n}\n"}
